# "1st testcase of forgotPassword scenario added"
#
# Test Suite sheet currently has:
#   Row1: TesCaseID | Description                     | Runmode
#   Row2: Registration | Registration suite description | YES
#   Row3: Login        | Catalogue suite description     | YES   <- wrong description (copy/paste leftover)
#
# This change:
#   - appends a new Row4 test case for the "Forgot Password" scenario
#   - fixes Row3's description to the correct "Login suite description"
#   - leaves the selection on C7 (where the user clicked next)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Forgot Password" suite test case row.
$ws.Range("A4").Value = "Forgot Password"

# Fix the mis-pasted description on the existing "Login" row.
$ws.Range("B3").Value = "Login suite description"

$ws.Range("B4").Value = "Forgot Password suite description"
$ws.Range("C4").Value = "YES"

# Leave the selection where the author left it.
[void]$ws.Range("C7").Select()
